$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 0.1012289004441671
$ws.Range("E2").Value = 6.81253326896009
$ws.Range("F2").Value = 17.16961920063468
